$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update name/lead text values in column L (rows 2-5)
$ws.Range("L2").Value = "dedic 6"
$ws.Range("L3").Value = "nadia 6"
$ws.Range("L4").Value = "cost 2"
$ws.Range("L5").Value = "jaenudin z"

# Update numeric values in column F (rows 4-5)
$ws.Range("F4").Value = 55210
$ws.Range("F5").Value = 55210
